$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "Pretrained Transformers as Universal Computation Engines"
$ws.Range("E28").Value = "https://ropiens.tistory.com/98"

$ws.Range("D32").Value = "정규 표현식 기초 (퍼옴)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/305"

$ws.Range("D37").Value = "[Paper Review]  Vq-wav2vec: Self-Supervised Learning of Discrete Speech Representations"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1452&mod=document&pageid=1"

$ws.Range("D39").Value = "A 3-Minute Review of PCA: Compression and Recovery"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/A-3-Minute-Review-of-PCA-Compression-and-Recovery-1"

$ws.Range("D51").Value = "[세이버메트릭스] 보살과 자살(척살)의 차이, 그리고 수비율"
$ws.Range("E51").Value = "https://bskyvision.com/1151"
